$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.189.83"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "'3.914.44"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'487.79"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("D6").Value = "'147.26"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.733"
$ws.Range("E9").Value = "  -0.28%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "'0.0000347"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "'42.97"
$ws.Range("E12").Value = "  -0.97%  "
$ws.Range("D13").Value = "'10.84"
$ws.Range("E13").Value = "  +4.59%  "
$ws.Range("D14").Value = "'4.547.21"
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "'3.924.56"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "'14.38"
$ws.Range("E16").Value = "  -5.78%  "
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").Value = "'19.89"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'1.13"
$ws.Range("E19").Value = "  -2.45%  "
$ws.Range("D20").Value = "'68.378.67"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "'441.93"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'15.13"
$ws.Range("E22").Value = "  +3.97%  "
$ws.Range("D23").Value = "'3.47"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("D24").Value = "'88.08"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "'11.40"
$ws.Range("E25").Value = "  +17.42%  "
$ws.Range("D26").Value = "'11.52"
$ws.Range("E26").Value = "  +13.00%  "
$ws.Range("D27").Value = "'3.65"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").Value = "'38.55"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").Value = "'723.77"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "'13.83"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("E32").Value = "  -0.56%  "
$ws.Range("E33").Value = "  +3.48%  "
$ws.Range("D34").Value = "'6.29"
$ws.Range("E34").Value = "  +17.61%  "
$ws.Range("D35").Value = "'42.19"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D38").Value = "'0.410"
$ws.Range("E38").Value = "  +22.10%  "
$ws.Range("E39").Value = "  -2.43%  "
$ws.Range("D42").Value = "'3.25"
$ws.Range("E42").Value = "  +5.95%  "
$ws.Range("D43").Value = "'0.0482"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("E45").Value = "  -0.32%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = "  +5.73%  "
$ws.Range("D48").Value = "'3.43"
$ws.Range("E48").Value = "  -2.61%  "
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("D50").Value = "'0.0₆0348"
$ws.Range("E50").Value = "  +32.48%  "
$ws.Range("D51").Value = "'146.06"
$ws.Range("E51").Value = "  -0.34%  "

# Row 36: PEPE -> OKB
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").Value = "'61.23"
$ws.Range("E36").Value = "  +5.73%  "

# Row 37: OKB -> PEPE
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "'0.0₃0864"
$ws.Range("E37").Value = "  +8.57%  "

# Row 40: Fetch.AI -> Dai
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41: Dai -> Fetch.AI
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "'2.98"
$ws.Range("E41").Value = "  +16.35%  "
